$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source export was refreshed: some accounts were removed, some added,
# a couple of balances were updated, and the rows were re-sorted by balance (desc).
# Rebuild the full data block (A2:C218) to match the refreshed export,
# then remove the extra now-unused row so the trailing blank/footer rows shift up.

$rowCount = 217

# Column A holds account numbers with leading zeros, so force text format
# on the whole column before writing values, otherwise Excel would coerce
# them to numbers and strip the leading zeros.
$ws.Range("A2:A" + ($rowCount + 1)).NumberFormat = "@"

$data = New-Object 'object[,]' 217,3
$data[0,0] = "003921139"; $data[0,1] = "GEISA"; $data[0,2] = 99452.54
$data[1,0] = "005002457"; $data[1,1] = "ROSANGELA"; $data[1,2] = 99007.59
$data[2,0] = "004222784"; $data[2,1] = "RAFAEL"; $data[2,2] = 89593.48
$data[3,0] = "004752534"; $data[3,1] = "CARLOS"; $data[3,2] = 65559.47
$data[4,0] = "005101676"; $data[4,1] = "ELENI"; $data[4,2] = 48319.7
$data[5,0] = "004452476"; $data[5,1] = "IVONE"; $data[5,2] = 47904.93
$data[6,0] = "004207955"; $data[6,1] = "SILVANIA"; $data[6,2] = 47633.14
$data[7,0] = "004399832"; $data[7,1] = "EULER"; $data[7,2] = 46914.63
$data[8,0] = "005105172"; $data[8,1] = "VALDIVINO"; $data[8,2] = 42029.42
$data[9,0] = "004575632"; $data[9,1] = "ADELE"; $data[9,2] = 23566.46
$data[10,0] = "004224011"; $data[10,1] = "THOMAS"; $data[10,2] = 21582.44
$data[11,0] = "004983395"; $data[11,1] = "MARCELO"; $data[11,2] = 18398.8
$data[12,0] = "004454365"; $data[12,1] = "RAFAEL"; $data[12,2] = 13735.23
$data[13,0] = "004946542"; $data[13,1] = "ALESSANDRA"; $data[13,2] = 9000.83
$data[14,0] = "004361159"; $data[14,1] = "HFR"; $data[14,2] = 5714.31
$data[15,0] = "004488571"; $data[15,1] = "CARLOS"; $data[15,2] = 1000
$data[16,0] = "004334158"; $data[16,1] = "LEONE"; $data[16,2] = 994.66
$data[17,0] = "004487016"; $data[17,1] = "ROGERIO"; $data[17,2] = 921.71
$data[18,0] = "004392159"; $data[18,1] = "RODRIGO"; $data[18,2] = 900.21
$data[19,0] = "004855570"; $data[19,1] = "LUISA"; $data[19,2] = 895.19
$data[20,0] = "004975924"; $data[20,1] = "SERGIO"; $data[20,2] = 874.97
$data[21,0] = "004216504"; $data[21,1] = "WANDER"; $data[21,2] = 850.24
$data[22,0] = "004563252"; $data[22,1] = "FERNANDO"; $data[22,2] = 837.81
$data[23,0] = "005245032"; $data[23,1] = "ROSA"; $data[23,2] = 824.46
$data[24,0] = "004322719"; $data[24,1] = "GISELA"; $data[24,2] = 814.1
$data[25,0] = "002064834"; $data[25,1] = "RAFAELA"; $data[25,2] = 813.54
$data[26,0] = "004381180"; $data[26,1] = "HFR"; $data[26,2] = 743.31
$data[27,0] = "004870019"; $data[27,1] = "MARIA"; $data[27,2] = 723.41
$data[28,0] = "002697806"; $data[28,1] = "CLAUDIA"; $data[28,2] = 705.53
$data[29,0] = "004388077"; $data[29,1] = "WLADMIR"; $data[29,2] = 703.17
$data[30,0] = "004359408"; $data[30,1] = "HEPTA"; $data[30,2] = 693.95
$data[31,0] = "004448303"; $data[31,1] = "NASSIM"; $data[31,2] = 692.74
$data[32,0] = "004972070"; $data[32,1] = "MARIA"; $data[32,2] = 670.51
$data[33,0] = "001882235"; $data[33,1] = "LAGO"; $data[33,2] = 661.4
$data[34,0] = "004252768"; $data[34,1] = "ALESSANDRO"; $data[34,2] = 656.22
$data[35,0] = "004481463"; $data[35,1] = "MARA"; $data[35,2] = 637.03
$data[36,0] = "004228456"; $data[36,1] = "FLASH"; $data[36,2] = 611.4
$data[37,0] = "004517080"; $data[37,1] = "TATIANA"; $data[37,2] = 607.94
$data[38,0] = "005079311"; $data[38,1] = "JOVINO"; $data[38,2] = 603.35
$data[39,0] = "004346716"; $data[39,1] = "TIAGO"; $data[39,2] = 598.88
$data[40,0] = "004813088"; $data[40,1] = "JULIANA"; $data[40,2] = 586.9
$data[41,0] = "004806244"; $data[41,1] = "CARLA"; $data[41,2] = 585.78
$data[42,0] = "005142592"; $data[42,1] = "ALBERTO"; $data[42,2] = 551.37
$data[43,0] = "004415557"; $data[43,1] = "FILIPE"; $data[43,2] = 545.71
$data[44,0] = "005055865"; $data[44,1] = "G3C"; $data[44,2] = 526.02
$data[45,0] = "004395314"; $data[45,1] = "MARIA"; $data[45,2] = 522.36
$data[46,0] = "005348975"; $data[46,1] = "JULIA"; $data[46,2] = 510
$data[47,0] = "004556853"; $data[47,1] = "MARCEL"; $data[47,2] = 483.1
$data[48,0] = "004387250"; $data[48,1] = "MONICA"; $data[48,2] = 478
$data[49,0] = "005338054"; $data[49,1] = "ELAINE"; $data[49,2] = 475.92
$data[50,0] = "000772433"; $data[50,1] = "MARCELO"; $data[50,2] = 465.51
$data[51,0] = "005266369"; $data[51,1] = "EG"; $data[51,2] = 459.39
$data[52,0] = "005203562"; $data[52,1] = "ROBERIO"; $data[52,2] = 453.66
$data[53,0] = "004398253"; $data[53,1] = "EULER"; $data[53,2] = 445.36
$data[54,0] = "004474776"; $data[54,1] = "GILSON"; $data[54,2] = 437
$data[55,0] = "004328934"; $data[55,1] = "VALERIA"; $data[55,2] = 424.53
$data[56,0] = "004556150"; $data[56,1] = "MARINA"; $data[56,2] = 409.35
$data[57,0] = "003435941"; $data[57,1] = "HEITOR"; $data[57,2] = 400
$data[58,0] = "004480970"; $data[58,1] = "ALBERTO"; $data[58,2] = 399.15
$data[59,0] = "005003629"; $data[59,1] = "ANDRE"; $data[59,2] = 370.24
$data[60,0] = "004587511"; $data[60,1] = "CARLOS"; $data[60,2] = 352.61
$data[61,0] = "004289402"; $data[61,1] = "LARISSA"; $data[61,2] = 349.86
$data[62,0] = "004381415"; $data[62,1] = "JOAO"; $data[62,2] = 349.74
$data[63,0] = "005009992"; $data[63,1] = "ALINE"; $data[63,2] = 330.17
$data[64,0] = "004214592"; $data[64,1] = "MERG"; $data[64,2] = 312.75
$data[65,0] = "004332103"; $data[65,1] = "JOSE"; $data[65,2] = 300.31
$data[66,0] = "005121919"; $data[66,1] = "JORGE"; $data[66,2] = 297.95
$data[67,0] = "004935287"; $data[67,1] = "ODILON"; $data[67,2] = 297.77
$data[68,0] = "004259659"; $data[68,1] = "BENTO"; $data[68,2] = 293.95
$data[69,0] = "004643737"; $data[69,1] = "LARA"; $data[69,2] = 289.88
$data[70,0] = "000330949"; $data[70,1] = "RENATO"; $data[70,2] = 285.37
$data[71,0] = "004641487"; $data[71,1] = "LAILA"; $data[71,2] = 283.14
$data[72,0] = "004451652"; $data[72,1] = "MATEUS"; $data[72,2] = 281.69
$data[73,0] = "000626491"; $data[73,1] = "FELIPE"; $data[73,2] = 280.07
$data[74,0] = "004355790"; $data[74,1] = "MINEIA"; $data[74,2] = 279.74
$data[75,0] = "004383190"; $data[75,1] = "MAFALDA"; $data[75,2] = 273.6
$data[76,0] = "004486497"; $data[76,1] = "ELENA"; $data[76,2] = 257.45
$data[77,0] = "004927044"; $data[77,1] = "CINTIA"; $data[77,2] = 257.37
$data[78,0] = "004580355"; $data[78,1] = "LARISSA"; $data[78,2] = 227.65
$data[79,0] = "004870976"; $data[79,1] = "HFR"; $data[79,2] = 222.74
$data[80,0] = "004884046"; $data[80,1] = "WILSON"; $data[80,2] = 218.7
$data[81,0] = "004526450"; $data[81,1] = "MSD"; $data[81,2] = 205.85
$data[82,0] = "004466221"; $data[82,1] = "WALTER"; $data[82,2] = 190.76
$data[83,0] = "004475395"; $data[83,1] = "DAVID"; $data[83,2] = 185.02
$data[84,0] = "005312963"; $data[84,1] = "ALAN"; $data[84,2] = 165.72
$data[85,0] = "004360431"; $data[85,1] = "CARLOS"; $data[85,2] = 164.01
$data[86,0] = "004213373"; $data[86,1] = "ALEXANDRE"; $data[86,2] = 162.27
$data[87,0] = "005044389"; $data[87,1] = "CLAUDIA"; $data[87,2] = 158.84
$data[88,0] = "004511696"; $data[88,1] = "KRYSCIA"; $data[88,2] = 150.47
$data[89,0] = "005274028"; $data[89,1] = "RAFAEL"; $data[89,2] = 149.15
$data[90,0] = "005305448"; $data[90,1] = "ALPHASITIO"; $data[90,2] = 139.53
$data[91,0] = "005141215"; $data[91,1] = "KARINA"; $data[91,2] = 137.66
$data[92,0] = "004243043"; $data[92,1] = "SUELI"; $data[92,2] = 134.8
$data[93,0] = "004237325"; $data[93,1] = "RICARDO"; $data[93,2] = 129.2
$data[94,0] = "004435987"; $data[94,1] = "MARCO"; $data[94,2] = 125.33
$data[95,0] = "004211911"; $data[95,1] = "ZENILDA"; $data[95,2] = 120
$data[96,0] = "004404342"; $data[96,1] = "ADSON"; $data[96,2] = 115.85
$data[97,0] = "004754920"; $data[97,1] = "LUIS"; $data[97,2] = 114.69
$data[98,0] = "004421636"; $data[98,1] = "PATRICIA"; $data[98,2] = 110
$data[99,0] = "004221638"; $data[99,1] = "CAROLINE"; $data[99,2] = 109.24
$data[100,0] = "004536602"; $data[100,1] = "TATIANY"; $data[100,2] = 108.62
$data[101,0] = "002687737"; $data[101,1] = "JOSE"; $data[101,2] = 101.13
$data[102,0] = "004547722"; $data[102,1] = "MARCIA"; $data[102,2] = 100
$data[103,0] = "004908680"; $data[103,1] = "ELENE"; $data[103,2] = 99.31
$data[104,0] = "004472076"; $data[104,1] = "RUBENS"; $data[104,2] = 99.18
$data[105,0] = "004339183"; $data[105,1] = "JALISON"; $data[105,2] = 95.69
$data[106,0] = "004431591"; $data[106,1] = "MARIO"; $data[106,2] = 93.87
$data[107,0] = "005256849"; $data[107,1] = "SANDRO"; $data[107,2] = 92.78
$data[108,0] = "004335031"; $data[108,1] = "EDMUNDO"; $data[108,2] = 92.73
$data[109,0] = "004350197"; $data[109,1] = "GISELA"; $data[109,2] = 91.94
$data[110,0] = "004239387"; $data[110,1] = "LUIZ"; $data[110,2] = 89.82
$data[111,0] = "004212132"; $data[111,1] = "JOAO"; $data[111,2] = 86.38
$data[112,0] = "004207374"; $data[112,1] = "ANGELICA"; $data[112,2] = 85.13
$data[113,0] = "005035754"; $data[113,1] = "JOSE"; $data[113,2] = 83.31
$data[114,0] = "004206790"; $data[114,1] = "EMMANUELLE"; $data[114,2] = 82.16
$data[115,0] = "004216657"; $data[115,1] = "JOAO"; $data[115,2] = 80.63
$data[116,0] = "004318604"; $data[116,1] = "RENAN"; $data[116,2] = 80.51
$data[117,0] = "004451996"; $data[117,1] = "ADRIANO"; $data[117,2] = 80.36
$data[118,0] = "004267976"; $data[118,1] = "E3"; $data[118,2] = 79.84
$data[119,0] = "001294033"; $data[119,1] = "VIVIANE"; $data[119,2] = 79.82
$data[120,0] = "005073033"; $data[120,1] = "NILBORN"; $data[120,2] = 79.39
$data[121,0] = "005009922"; $data[121,1] = "ANA"; $data[121,2] = 79.02
$data[122,0] = "004470679"; $data[122,1] = "RODOLFO"; $data[122,2] = 77.51
$data[123,0] = "004565108"; $data[123,1] = "GUSTAVO"; $data[123,2] = 75.18
$data[124,0] = "005133039"; $data[124,1] = "PAULO"; $data[124,2] = 66.51
$data[125,0] = "004855596"; $data[125,1] = "MARIANA"; $data[125,2] = 64.36
$data[126,0] = "004335251"; $data[126,1] = "EDMUNDO"; $data[126,2] = 62.39
$data[127,0] = "000834301"; $data[127,1] = "MARCUS"; $data[127,2] = 57.13
$data[128,0] = "004588677"; $data[128,1] = "RACHEL"; $data[128,2] = 55.91
$data[129,0] = "004517506"; $data[129,1] = "LUIZ"; $data[129,2] = 55.87
$data[130,0] = "004215217"; $data[130,1] = "CAROLINA"; $data[130,2] = 55.66
$data[131,0] = "004321092"; $data[131,1] = "DANIEL"; $data[131,2] = 55.23
$data[132,0] = "004329229"; $data[132,1] = "GABRIEL"; $data[132,2] = 54.91
$data[133,0] = "004999434"; $data[133,1] = "EDUARDO"; $data[133,2] = 54.21
$data[134,0] = "005032151"; $data[134,1] = "ANA"; $data[134,2] = 52.9
$data[135,0] = "004268684"; $data[135,1] = "PATRICIA"; $data[135,2] = 52.7
$data[136,0] = "002277249"; $data[136,1] = "DANILO"; $data[136,2] = 52.44
$data[137,0] = "004400640"; $data[137,1] = "FELIPE"; $data[137,2] = 51.44
$data[138,0] = "004115403"; $data[138,1] = "HEBERT"; $data[138,2] = 50.87
$data[139,0] = "004208447"; $data[139,1] = "LEILA"; $data[139,2] = 50
$data[140,0] = "004278033"; $data[140,1] = "DAISY"; $data[140,2] = 47.37
$data[141,0] = "004491730"; $data[141,1] = "DENISE"; $data[141,2] = 47.19
$data[142,0] = "005216881"; $data[142,1] = "RENAN"; $data[142,2] = 46.76
$data[143,0] = "001719494"; $data[143,1] = "LUIS"; $data[143,2] = 46.74
$data[144,0] = "004277637"; $data[144,1] = "LARA"; $data[144,2] = 46.2
$data[145,0] = "001731007"; $data[145,1] = "GUILHERME"; $data[145,2] = 44.59
$data[146,0] = "004581652"; $data[146,1] = "CINCO"; $data[146,2] = 44.13
$data[147,0] = "004805133"; $data[147,1] = "PATRICIA"; $data[147,2] = 41.48
$data[148,0] = "004958578"; $data[148,1] = "ASSAKO"; $data[148,2] = 40.7
$data[149,0] = "005165116"; $data[149,1] = "ANA"; $data[149,2] = 40.11
$data[150,0] = "004998717"; $data[150,1] = "GIOVANE"; $data[150,2] = 40.08
$data[151,0] = "004238164"; $data[151,1] = "DANIELA"; $data[151,2] = 38.3
$data[152,0] = "004520100"; $data[152,1] = "ALEXANDRE"; $data[152,2] = 37.96
$data[153,0] = "002401479"; $data[153,1] = "JULIO"; $data[153,2] = 37.84
$data[154,0] = "005000656"; $data[154,1] = "LUCIA"; $data[154,2] = 35.88
$data[155,0] = "004211922"; $data[155,1] = "CARLOS"; $data[155,2] = 34.71
$data[156,0] = "004340984"; $data[156,1] = "RENATA"; $data[156,2] = 34
$data[157,0] = "004994036"; $data[157,1] = "BALTASAR"; $data[157,2] = 33.73
$data[158,0] = "004691225"; $data[158,1] = "ANNA"; $data[158,2] = 33.64
$data[159,0] = "004472431"; $data[159,1] = "LUIS"; $data[159,2] = 33.08
$data[160,0] = "005018038"; $data[160,1] = "ELAINE"; $data[160,2] = 27.91
$data[161,0] = "004377415"; $data[161,1] = "ANGELA"; $data[161,2] = 26.37
$data[162,0] = "004240292"; $data[162,1] = "MARCO"; $data[162,2] = 24.3
$data[163,0] = "004404724"; $data[163,1] = "LEANDRO"; $data[163,2] = 24.14
$data[164,0] = "005173958"; $data[164,1] = "VENIA"; $data[164,2] = 23.27
$data[165,0] = "005186167"; $data[165,1] = "ANDREA"; $data[165,2] = 21
$data[166,0] = "004214604"; $data[166,1] = "MARIA"; $data[166,2] = 20.75
$data[167,0] = "004920447"; $data[167,1] = "MARILIA"; $data[167,2] = 20.67
$data[168,0] = "004458604"; $data[168,1] = "FABIOLA"; $data[168,2] = 20
$data[169,0] = "004204255"; $data[169,1] = "AMADO"; $data[169,2] = 18.77
$data[170,0] = "004368994"; $data[170,1] = "CRISTINA"; $data[170,2] = 18.56
$data[171,0] = "004756968"; $data[171,1] = "DANIELY"; $data[171,2] = 18.08
$data[172,0] = "001879977"; $data[172,1] = "THAISSA"; $data[172,2] = 17.14
$data[173,0] = "005143579"; $data[173,1] = "GABRIEL"; $data[173,2] = 16.18
$data[174,0] = "005169333"; $data[174,1] = "EDUARDO"; $data[174,2] = 16.12
$data[175,0] = "004422594"; $data[175,1] = "WANDIR"; $data[175,2] = 14.67
$data[176,0] = "000827730"; $data[176,1] = "LUCIANA"; $data[176,2] = 13.29
$data[177,0] = "004752461"; $data[177,1] = "SERGIO"; $data[177,2] = 10.77
$data[178,0] = "004216298"; $data[178,1] = "FLORDELIZ"; $data[178,2] = 9.75
$data[179,0] = "004527606"; $data[179,1] = "MARCIA"; $data[179,2] = 9.52
$data[180,0] = "004264780"; $data[180,1] = "MARCELO"; $data[180,2] = 8.99
$data[181,0] = "004921978"; $data[181,1] = "ELAINE"; $data[181,2] = 8.08
$data[182,0] = "004214460"; $data[182,1] = "MARIA"; $data[182,2] = 7.54
$data[183,0] = "004693631"; $data[183,1] = "NELY"; $data[183,2] = 7.36
$data[184,0] = "004530494"; $data[184,1] = "ROSANGELA"; $data[184,2] = 6.94
$data[185,0] = "004854496"; $data[185,1] = "JOSE"; $data[185,2] = 6.64
$data[186,0] = "004448501"; $data[186,1] = "JOAO"; $data[186,2] = 5.55
$data[187,0] = "004216434"; $data[187,1] = "JAIME"; $data[187,2] = 4.83
$data[188,0] = "005142624"; $data[188,1] = "RODRIGO"; $data[188,2] = 4.75
$data[189,0] = "004239624"; $data[189,1] = "NINA"; $data[189,2] = 4.29
$data[190,0] = "004848927"; $data[190,1] = "ULDARICO"; $data[190,2] = 3.62
$data[191,0] = "005142661"; $data[191,1] = "SABRINA"; $data[191,2] = 3.6
$data[192,0] = "004382374"; $data[192,1] = "THEOMAR"; $data[192,2] = 3.41
$data[193,0] = "004335144"; $data[193,1] = "EDMUNDO"; $data[193,2] = 2.72
$data[194,0] = "005341184"; $data[194,1] = "BRENO"; $data[194,2] = 1.85
$data[195,0] = "004332783"; $data[195,1] = "IRON"; $data[195,2] = 1.73
$data[196,0] = "004886366"; $data[196,1] = "RENATO"; $data[196,2] = 1.57
$data[197,0] = "005366255"; $data[197,1] = "RAPHAELA"; $data[197,2] = 1.39
$data[198,0] = "004308815"; $data[198,1] = "ZELI"; $data[198,2] = 1.25
$data[199,0] = "005228239"; $data[199,1] = "DEBORA"; $data[199,2] = 0.85
$data[200,0] = "004223502"; $data[200,1] = "BRUNA"; $data[200,2] = 0.78
$data[201,0] = "004212581"; $data[201,1] = "MARIA"; $data[201,2] = 0.59
$data[202,0] = "004550605"; $data[202,1] = "REJANE"; $data[202,2] = 0.53
$data[203,0] = "004453302"; $data[203,1] = "ISABELLA"; $data[203,2] = 0.39
$data[204,0] = "004806286"; $data[204,1] = "VERA"; $data[204,2] = 0.19
$data[205,0] = "004371857"; $data[205,1] = "NAZARETH"; $data[205,2] = 0.18
$data[206,0] = "004332207"; $data[206,1] = "IRACY"; $data[206,2] = 0.16
$data[207,0] = "004357159"; $data[207,1] = "JOAO"; $data[207,2] = 0.15
$data[208,0] = "004320840"; $data[208,1] = "NATALIA"; $data[208,2] = 0.14
$data[209,0] = "004466350"; $data[209,1] = "RAQUEL"; $data[209,2] = 0.11
$data[210,0] = "005047946"; $data[210,1] = "GABRIEL"; $data[210,2] = 0.09
$data[211,0] = "004589311"; $data[211,1] = "CLARICE"; $data[211,2] = 0.06
$data[212,0] = "004321016"; $data[212,1] = "JOAQUIM"; $data[212,2] = 0.02
$data[213,0] = "004850070"; $data[213,1] = "RENATO"; $data[213,2] = 0.02
$data[214,0] = "002878817"; $data[214,1] = "GUILHERME"; $data[214,2] = 0.01
$data[215,0] = "004400000"; $data[215,1] = "VILMA"; $data[215,2] = 0.01
$data[216,0] = "004612043"; $data[216,1] = "YURI"; $data[216,2] = 0.01

$ws.Range("A2:C" + ($rowCount + 1)).Value = $data

# The refreshed export has one fewer data row than before, so delete the
# row right after the new data block to shift the trailing blank row and
# the "Filtros aplicados" footer row back up into place.
$ws.Rows.Item($rowCount + 2).Delete()

Write-Host "Rewrote" $rowCount "data rows."
